$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (shifts existing rows 6-11 down to 7-12)
$ws.Rows(6).Insert()

# Populate the new row 6 with the "verify emailid" test step
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "verify emailid"
$ws.Range("C6").Value = "input"
$ws.Range("E6").Value = "email"
$ws.Range("F6").Value = "yes"

# Flip column F (runmode) from "yes" to "no" for the affected rows,
# except row 9 ("click on send button"), which stays "yes"
$ws.Range("F5").Value = "no"
$ws.Range("F7").Value = "no"
$ws.Range("F8").Value = "no"
$ws.Range("F10").Value = "no"
$ws.Range("F11").Value = "no"

# Set the new email value last, so shared-string ordering matches
$ws.Range("D6").Value = "abc.xyz04071991@gmail.com"

# Add hyperlink for the new email address in D6
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:abc.xyz04071991@gmail.com")

# Update the active selection to F9
$ws.Range("F9").Select()
